# Petty cash book update — 4 Mar 2021, end of day.
# Sheet1 ("Buku KAS HARIAN" petty-cash ledger) rows 20-24 are corrected:
#  - row 20: the stray text placeholder in C20 is removed and the actual
#    debit amount (10,000) is recorded in D20.
#  - row 21: unchanged inputs, formula now recalculates cleanly.
#  - row 22: a Wages Expense debit of 45,000 is recorded.
#  - row 23: a new TRANSFER BCA debit (450000+487500) is recorded.
#  - row 24: a new FREIGHT OUT debit of 65,000 is recorded.
# Fixing the broken C20 text entry lets the whole E-column running-balance
# formula chain (rows 20-113) recalculate to real numbers instead of
# #VALUE! errors.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 20: clear the stray text in C20, fix the debit amount in D20.
$ws.Range("C20").Clear()
$ws.Range("D20").Value = 10000

# Row 23: new TRANSFER BCA entry.
$ws.Range("B23").Value = "TRANSFER BCA"
$ws.Range("D23").Formula = "=450000+487500"

# Row 24: new FREIGHT OUT entry.
$ws.Range("B24").Value = "FREIGHT OUT"
$ws.Range("D24").Formula = "=65000"

# Row 22: Wages Expense debit.
$ws.Range("D22").Formula = "=45000"

# Re-touch the shared running-balance formulas for the rows whose inputs
# just changed (E20:E24), so each recalculates against the corrected
# figures above instead of keeping a stale cached #VALUE! result. Once
# E24 is fresh, the rest of the E-column chain (rows 25-113) ripples
# through correctly on its own.
$ws.Range("E20").Formula = "=E19+C20-D20"
$ws.Range("E21").Formula = "=E20+C21-D21"
$ws.Range("E22").Formula = "=E21+C22-D22"
$ws.Range("E23").Formula = "=E22+C23-D23"
$ws.Range("E24").Formula = "=E23+C24-D24"

# Update the remembered selection to match the author's last position.
$ws.Range("D21").Select()
